# Add the two missing log entries to the "宋明硕" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("宋明硕")

$ws.Range("A21").Value = "2019/5/31 18:00-20:00"
$ws.Range("B21").Value = "简单项目框架"
$ws.Range("A22").Value = "2019/6/11 22:00-24:00"
$ws.Range("B22").Value = "提交文档"

# Mirror the selection change recorded for this sheet after the edit.
$ws.Range("A23").Select()
